# Update computed price/profit columns (H:N) on each Leve table per the
# scheduled-runner refresh. Values come straight from the source data feed;
# there are no formulas in these sheets, so we just overwrite the cells.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 2368.6562
$ws.Range("I80").Value = 802.7778
$ws.Range("J80").Value = 2981.3914
$ws.Range("K80").Value = 2408.3334
$ws.Range("L80").Value = 8944.174199999999
$ws.Range("M80").Value = -1410.3334
$ws.Range("N80").Value = -10940.1742
# Row 83
$ws.Range("H83").Value = 2368.6562
$ws.Range("I83").Value = 802.7778
$ws.Range("J83").Value = 2981.3914
$ws.Range("K83").Value = 7225.000199999999
$ws.Range("L83").Value = 26832.5226
$ws.Range("M83").Value = -2233.000199999999
$ws.Range("N83").Value = -36816.5226
# Row 100
$ws.Range("H100").Value = 5150.1514
$ws.Range("I100").Value = 3781.5557
$ws.Range("K100").Value = 3781.5557
$ws.Range("M100").Value = -3240.5557
# Row 107
$ws.Range("H107").Value = 1404.4615
$ws.Range("J107").Value = 1507.5
$ws.Range("L107").Value = 1507.5
$ws.Range("N107").Value = -5347.5
# Row 112
$ws.Range("H112").Value = 1715.5186
$ws.Range("I112").Value = 1056
$ws.Range("J112").Value = 1865.409
$ws.Range("K112").Value = 3168
$ws.Range("L112").Value = 5596.227000000001
$ws.Range("M112").Value = -2060
$ws.Range("N112").Value = -7812.227000000001
# Row 125
$ws.Range("H125").Value = 1169.125
$ws.Range("I125").Value = 1277.125
$ws.Range("J125").Value = 1061.125
$ws.Range("K125").Value = 11494.125
$ws.Range("L125").Value = 9550.125
$ws.Range("M125").Value = -9034.125
$ws.Range("N125").Value = -14470.125
# Row 137
$ws.Range("H137").Value = 6499.15
$ws.Range("I137").Value = 8019.6665
$ws.Range("J137").Value = 4218.375
$ws.Range("K137").Value = 24058.9995
$ws.Range("L137").Value = 12655.125
$ws.Range("M137").Value = -21508.9995
$ws.Range("N137").Value = -17755.125
# Row 138
$ws.Range("H138").Value = 2682.9656
$ws.Range("I138").Value = 1550.6471
$ws.Range("J138").Value = 4287.0835
$ws.Range("K138").Value = 4651.9413
$ws.Range("L138").Value = 12861.2505
$ws.Range("M138").Value = 488.0587000000005
$ws.Range("N138").Value = -23141.2505

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3546.5386
$ws.Range("I61").Value = 1620
$ws.Range("J61").Value = 5197.857
$ws.Range("K61").Value = 1620
$ws.Range("L61").Value = 5197.857
$ws.Range("M61").Value = -1408
$ws.Range("N61").Value = -5621.857
# Row 74
$ws.Range("H74").Value = 1219.3334
$ws.Range("I74").Value = 1219.3334
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1219.3334
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -345.3334
$ws.Range("N74").Value = $null
# Row 77
$ws.Range("H77").Value = 1219.3334
$ws.Range("I77").Value = 1219.3334
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 6096.666999999999
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -1728.666999999999
$ws.Range("N77").Value = $null
# Row 110
$ws.Range("H110").Value = 5154.727
$ws.Range("I110").Value = 5226.5
$ws.Range("K110").Value = 5226.5
$ws.Range("M110").Value = -3181.5
# Row 132
$ws.Range("H132").Value = 2239.6155
$ws.Range("I132").Value = 2009.5834
$ws.Range("K132").Value = 6028.7502
$ws.Range("M132").Value = -3498.7502
# Row 136
$ws.Range("H136").Value = 3546.5386
$ws.Range("I136").Value = 1620
$ws.Range("J136").Value = 5197.857
$ws.Range("K136").Value = 4860
$ws.Range("L136").Value = 15593.571
$ws.Range("M136").Value = -2310
$ws.Range("N136").Value = -20693.571

$ws = $wb.Worksheets.Item("BSM")
# Row 110
$ws.Range("H110").Value = 110000
$ws.Range("J110").Value = 110000
$ws.Range("L110").Value = 110000
$ws.Range("N110").Value = -118180

$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Range("H132").Value = 3661.7778
$ws.Range("I132").Value = 3661.7778
$ws.Range("K132").Value = 10985.3334
$ws.Range("M132").Value = -8455.3334

$ws = $wb.Worksheets.Item("CUL")
# Row 13
$ws.Range("H13").Value = 987.25
$ws.Range("I13").Value = 987.25
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 2961.75
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -2793.75
$ws.Range("N13").Value = $null
# Row 39
$ws.Range("H39").Value = 1248.6
$ws.Range("J39").Value = 1417.2
$ws.Range("L39").Value = 4251.6
$ws.Range("N39").Value = -4839.6

$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 6285.294
$ws.Range("I97").Value = 457.36365
$ws.Range("K97").Value = 457.36365
$ws.Range("M97").Value = 38.63634999999999
# Row 122
$ws.Range("H122").Value = 6071.143
$ws.Range("I122").Value = 5749.5
$ws.Range("K122").Value = 17248.5
$ws.Range("M122").Value = -14798.5
# Row 132
$ws.Range("H132").Value = 8216.333000000001
$ws.Range("I132").Value = 7769.885
$ws.Range("K132").Value = 23309.655
$ws.Range("M132").Value = -20779.655

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 1446.75
$ws.Range("I61").Value = 1496.2858
$ws.Range("K61").Value = 1496.2858
$ws.Range("M61").Value = -1294.2858
# Row 93
$ws.Range("H93").Value = 2498.4
$ws.Range("I93").Value = 2197.6667
$ws.Range("J93").Value = 2949.5
$ws.Range("K93").Value = 2197.6667
$ws.Range("L93").Value = 2949.5
$ws.Range("M93").Value = -949.6667000000002
$ws.Range("N93").Value = -5445.5
# Row 113
$ws.Range("H113").Value = 1446.75
$ws.Range("I113").Value = 1496.2858
$ws.Range("K113").Value = 1496.2858
$ws.Range("M113").Value = 673.7141999999999
# Row 136
$ws.Range("H136").Value = 2613.8147
$ws.Range("J136").Value = 3801.4443
$ws.Range("L136").Value = 11404.3329
$ws.Range("N136").Value = -16504.3329

$ws = $wb.Worksheets.Item("WVR")
# Row 41
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = $null
# Row 46
$ws.Range("H46").Value = 65457
$ws.Range("J46").Value = 65457
$ws.Range("L46").Value = 65457
$ws.Range("N46").Value = -65919
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = $null
# Row 132
$ws.Range("H132").Value = 2180.4644
$ws.Range("I132").Value = 2076.037
$ws.Range("K132").Value = 6228.110999999999
$ws.Range("M132").Value = -3698.110999999999
# Row 134
$ws.Range("H134").Value = 65457
$ws.Range("J134").Value = 65457
$ws.Range("L134").Value = 196371
$ws.Range("N134").Value = -201441
# Row 136
$ws.Range("H136").Value = 75013.53999999999
$ws.Range("I136").Value = 75013.53999999999
$ws.Range("K136").Value = 225040.62
$ws.Range("M136").Value = -222490.62
